$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.456.72"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "'2.982.69"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'381.62"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "'104.19"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'36.62"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "'3.455.22"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "'18.43"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "'7.79"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("D16").Value = "'2.983.19"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "'11.16"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "'0.996"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'51.472.93"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "'3.11"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "'12.58"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "'70.25"
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("D24").Value = "'267.04"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "'3.22"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").Value = "'7.82"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("D27").Value = "'7.30"
$ws.Range("E27").Value = "  -3.81%  "
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'26.06"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'10.38"
$ws.Range("E32").Value = "  +3.68%  "
$ws.Range("D33").Value = "'34.72"
$ws.Range("E33").Value = "  +3.77%  "
$ws.Range("D34").Value = "'51.40"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "'2.06"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").Value = "'0.0446"
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'3.29"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").Value = "'16.93"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.57"
$ws.Range("E40").Value = "  +4.06%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.117"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.84"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "'3.83"
$ws.Range("E43").Value = "  +12.54%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'122.25"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Value = "'0.271"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").Value = "'2.024.16"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("D50").Value = "'3.280.17"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("E51").Value = "  +0.33%  "
